$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2,3) {
    $ws.Range("D$r").Value = 0.0762
    $ws.Range("E$r").Value = 0.131
    $ws.Range("G$r").Value = 0.0712707182320442
    $ws.Range("H$r").Value = 0.0712707182320442
    $ws.Range("I$r").Value = 0.08176795580110496
    $ws.Range("J$r").Value = 0.07158867967076334
    $ws.Range("K$r").Value = 1.28
    $ws.Range("L$r").Value = 0.07071823204419889
    $ws.Range("M$r").Value = 0.915
    $ws.Range("N$r").Value = 0.06727941176470589
    $ws.Range("O$r").Value = 0.71484375
    $ws.Range("P$r").Value = 0.915
    $ws.Range("Q$r").Value = 0.06727941176470589
    $ws.Range("R$r").Value = 0.71484375
    $ws.Range("U$r").Value = 1.58
    $ws.Range("V$r").Value = 0.1161764705882353
    $ws.Range("W$r").Value = 1.174311926605504
    $ws.Range("X$r").Value = 0.05529453638342827
    $ws.Range("Y$r").Value = 1.119017390222076
    $ws.Range("Z$r").Value = 36.2
    $ws.Range("AA$r").Value = 2.591510204081633
    $ws.Range("AB$r").Value = 0.05487765993591989
    $ws.Range("AC$r").Value = 2.536632544145713
    $ws.Range("AD$r").Value = 0.253
    $ws.Range("AF$r").Value = 0.253
    $ws.Range("AG$r").Value = -1.327
    $ws.Range("AH$r").Value = 0.01826319208835631
    $ws.Range("AI$r").Value = 0.1435053885422575
    $ws.Range("AJ$r").Value = -0.1081235231809664
    $ws.Range("AK$r").Value = -7.251366120218577
    $ws.Range("AL$r").Value = 0.014
    $ws.Range("AM$r").Value = 0.006
    $ws.Range("AN$r").Value = 0.1675496688741722
    $ws.Range("AO$r").Value = 105.7142857142857
    $ws.Range("AP$r").Value = -0.8788079470198675
    $ws.Range("AQ$r").Value = 246.6666666666667
}
